$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.717000000000001
$ws.Range("E4").Value = 13.421

$ws.Range("E5").Value = 13.41

$ws.Range("D6").Value = -8.264000000000001

$ws.Range("D7").Value = -7.523999999999999

$ws.Range("E8").Value = 13.718

$ws.Range("D16").Value = -8.045999999999999
$ws.Range("E16").Value = 12.978

$ws.Range("D20").Value = -8.184999999999999

$ws.Range("E22").Value = 13.495
